{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Para[0] = intro sentence, Para[1] = empty paragraph (to receive new text),\n// Para[2] = paragraph holding the _GoBack bookmark.\nconst introPara = paragraphs.items[0];\n\n// Insert a brand-new paragraph containing the new sentence right after the\n// intro paragraph (this mirrors how Word materializes a freshly typed\n// paragraph, including its run-level language formatting).\nintroPara.insertParagraph(\"A new sentence in this file.\", Word.InsertLocation.after);\nawait context.sync();\n\n// The previously-empty second paragraph is now redundant (we added a new\n// paragraph instead of typing into it) -- remove it so the bookmark\n// paragraph directly follows our new sentence, matching the target layout.\nconst paragraphsAfterInsert = body.paragraphs;\nparagraphsAfterInsert.load(\"items\");\nawait context.sync();\nparagraphsAfterInsert.items[2].delete();\nawait context.sync();\n\n// Finally, append two new empty paragraphs at the very end of the document\n// body (after the bookmark paragraph, before the section break).\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 1 = the intro sentence (\"XTM Change Control Connector ...\").\n# Paragraph 2 = the (currently empty) paragraph that is getting the new\n# sentence. Paragraph 3 = the paragraph holding the _GoBack bookmark.\n$introPara = $d.Paragraphs.Item(1)\n\n# Insert a brand-new paragraph right after the intro paragraph and give it\n# the new sentence's text. This mirrors how Word materializes a freshly\n# typed paragraph (including the inherited run-level language formatting).\n$introPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item(2)\n$newPara.Range.Text = \"A new sentence in this file.\"\n\n# The original empty paragraph (now pushed down to index 3) is redundant,\n# so remove it -- the bookmark paragraph should directly follow our new\n# sentence paragraph.\n$d.Paragraphs.Item(3).Range.Delete()\n\n# Finally, append two new empty paragraphs at the very end of the document\n# (after the bookmark paragraph, before the section break).\n$endRange = $d.Range()\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n\n$endRange2 = $d.Range()\n$endRange2.Collapse(0)\n$endRange2.InsertParagraphAfter()\n"}
